$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("slots")
$ws.Range("B2").Value = "favorite-sandwich"
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
